{"js": "  // Replacement pairs extracted from the diff: [oldText, newText], in document order.\n  const pairs = [\n    [\"32+48=80\", \"7+58=65\"],\n    [\"26-13=13\", \"34-30=4\"],\n    [\"69+4=73\", \"48+6=54\"],\n    [\"45-1=44\", \"5+44=49\"],\n    [\"30+52=82\", \"69+3=72\"],\n    [\"98-3=95\", \"43+52=95\"],\n    [\"64-41=23\", \"5+37=42\"],\n    [\"82-42=40\", \"94-73=21\"],\n    [\"26+17=43\", \"51-5=46\"],\n    [\"80+3=83\", \"8+42=50\"],\n    [\"58+32=90\", \"29+24=53\"],\n    [\"46+16=62\", \"22+66=88\"],\n    [\"35+58=93\", \"83-32=51\"],\n    [\"72+14=86\", \"73-0=73\"],\n    [\"15+7=22\", \"50-50=0\"],\n    [\"89-68=21\", \"98-84=14\"],\n    [\"18+72=90\", \"90-8=82\"],\n    [\"25+10=35\", \"16-8=8\"],\n    [\"34+4=38\", \"24+5=29\"],\n    [\"43-26=17\", \"54+18=72\"],\n    [\"16+19=35\", \"50-6=44\"],\n    [\"85-84=1\", \"64+22=86\"],\n    [\"84-23=61\", \"25+60=85\"],\n    [\"18+60=78\", \"29+56=85\"],\n    [\"4+62=66\", \"64-26=38\"],\n    [\"99-56=43\", \"40+18=58\"],\n    [\"8+56=64\", \"66-2=64\"],\n    [\"47-42=5\", \"90-24=66\"],\n    [\"89-17=72\", \"85-19=66\"],\n    [\"2+16=18\", \"29+8=37\"],\n    [\"34+26=60\", \"19-11=8\"],\n    [\"20+22=42\", \"48+0=48\"],\n    [\"26+57=83\", \"40+9=49\"],\n    [\"39-14=25\", \"98-11=87\"],\n    [\"98-59=39\", \"19+36=55\"],\n    [\"78-12=66\", \"80-25=55\"],\n    [\"6+66=72\", \"61+11=72\"],\n    [\"58+21=79\", \"98-81=17\"],\n    [\"83-23=60\", \"13+71=84\"],\n    [\"78-56=22\", \"77-21=56\"],\n    [\"62-20=42\", \"24+33=57\"],\n    [\"80+4=84\", \"24-0=24\"],\n    [\"35-0=35\", \"11-3=8\"],\n    [\"59-51=8\", \"0+52=52\"],\n    [\"92-41=51\", \"66-65=1\"],\n    [\"56-4=52\", \"94-48=46\"],\n    [\"22+72=94\", \"33+52=85\"],\n    [\"59-4=55\", \"73+2=75\"],\n    [\"98-86=12\", \"95+3=98\"],\n    [\"76-22=54\", \"10-7=3\"],\n    [\"94+3=97\", \"43+50=93\"],\n    [\"32+38=70\", \"59-15=44\"],\n    [\"94-23=71\", \"21+29=50\"],\n    [\"74-23=51\", \"42+30=72\"],\n    [\"81-20=61\", \"21+39=60\"],\n    [\"81-4=77\", \"31+23=54\"],\n    [\"57-17=40\", \"89-85=4\"],\n    [\"55-18=37\", \"44+54=98\"],\n    [\"15+80=95\", \"39-2=37\"],\n    [\"28-6=22\", \"86-22=64\"],\n    [\"38+47=85\", \"37+29=66\"],\n    [\"32+22=54\", \"44+46=90\"],\n    [\"23+75=98\", \"25-19=6\"],\n    [\"12+16=28\", \"82-63=19\"],\n    [\"28+29=57\", \"84-74=10\"],\n    [\"65-57=8\", \"28+45=73\"],\n    [\"54-10=44\", \"6+14=20\"],\n    [\"2+42=44\", \"50+21=71\"],\n    [\"37+30=67\", \"26+1=27\"],\n    [\"37+48=85\", \"24+48=72\"],\n    [\"16+43=59\", \"92-88=4\"],\n    [\"95-71=24\", \"84-81=3\"],\n    [\"44+14=58\", \"56-16=40\"],\n    [\"18-5=13\", \"61-40=21\"],\n    [\"27-2=25\", \"9+1=10\"],\n    [\"60-38=22\", \"12+58=70\"],\n    [\"89-47=42\", \"90-46=44\"],\n    [\"22-15=7\", \"52-46=6\"],\n    [\"91-0=91\", \"67-2=65\"],\n    [\"78-36=42\", \"78-64=14\"],\n    [\"61+28=89\", \"52+15=67\"],\n    [\"63-29=34\", \"40+44=84\"],\n    [\"70-43=27\", \"37-22=15\"],\n    [\"79-4=75\", \"89-80=9\"],\n    [\"9+57=66\", \"7+78=85\"],\n    [\"29+13=42\", \"0+96=96\"],\n    [\"80-20=60\", \"28+61=89\"],\n    [\"38-17=21\", \"26+25=51\"],\n    [\"81-43=38\", \"89-33=56\"],\n    [\"8+37=45\", \"16+12=28\"],\n    [\"19+30=49\", \"77+3=80\"],\n    [\"81+14=95\", \"7+11=18\"],\n    [\"79-67=12\", \"5+40=45\"],\n    [\"61-0=61\", \"16+53=69\"],\n    [\"3+16=19\", \"71-69=2\"],\n    [\"35+48=83\", \"34+49=83\"],\n    [\"44+53=97\", \"19+51=70\"],\n    [\"8+11=19\", \"41-15=26\"],\n    [\"33-2=31\", \"20+19=39\"],\n    [\"82+6=88\", \"21+31=52\"]\n  ];\n\n  const tables = context.document.body.tables;\n  tables.load(\"items\");\n  await context.sync();\n\n  if (tables.items.length === 0) {\n    throw new Error(\"No tables found in document body.\");\n  }\n\n  const table = tables.items[0];\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  // Load all cells of all rows up front.\n  const allCells = [];\n  for (const row of rows.items) {\n    row.cells.load(\"items\");\n  }\n  await context.sync();\n\n  for (const row of rows.items) {\n    for (const cell of row.cells.items) {\n      allCells.push(cell);\n    }\n  }\n\n  if (allCells.length !== pairs.length) {\n    throw new Error(\n      \"Cell count (\" + allCells.length + \") does not match expected pair count (\" + pairs.length + \").\"\n    );\n  }\n\n  // Grab the (single) paragraph of each cell and load its text so we can\n  // verify the expected \"before\" content prior to mutating it.\n  const paragraphs = allCells.map((cell) => cell.body.paragraphs.getFirst());\n  for (const paragraph of paragraphs) {\n    paragraph.load(\"text\");\n  }\n  await context.sync();\n\n  for (let i = 0; i < paragraphs.length; i++) {\n    const paragraph = paragraphs[i];\n    const [oldText, newText] = pairs[i];\n    const actualText = paragraph.text.trim();\n    if (actualText !== oldText) {\n      throw new Error(\n        \"Cell \" + i + \" text mismatch: expected '\" + oldText + \"' but found '\" + actualText + \"'.\"\n      );\n    }\n    // Replace the run's text in place (via its Range) so the existing\n    // run/paragraph formatting (font, size, alignment, ...) is preserved -\n    // this mirrors typing over the selected text in Word, rather than\n    // clearing the paragraph and inserting a brand-new, unformatted run.\n    const range = paragraph.getRange();\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n\n  await context.sync();\n", "ps1": "# Replacement pairs extracted from the diff: @(oldText, newText), in document order.\n$pairs = @(\n    @(\"32+48=80\", \"7+58=65\"),\n    @(\"26-13=13\", \"34-30=4\"),\n    @(\"69+4=73\", \"48+6=54\"),\n    @(\"45-1=44\", \"5+44=49\"),\n    @(\"30+52=82\", \"69+3=72\"),\n    @(\"98-3=95\", \"43+52=95\"),\n    @(\"64-41=23\", \"5+37=42\"),\n    @(\"82-42=40\", \"94-73=21\"),\n    @(\"26+17=43\", \"51-5=46\"),\n    @(\"80+3=83\", \"8+42=50\"),\n    @(\"58+32=90\", \"29+24=53\"),\n    @(\"46+16=62\", \"22+66=88\"),\n    @(\"35+58=93\", \"83-32=51\"),\n    @(\"72+14=86\", \"73-0=73\"),\n    @(\"15+7=22\", \"50-50=0\"),\n    @(\"89-68=21\", \"98-84=14\"),\n    @(\"18+72=90\", \"90-8=82\"),\n    @(\"25+10=35\", \"16-8=8\"),\n    @(\"34+4=38\", \"24+5=29\"),\n    @(\"43-26=17\", \"54+18=72\"),\n    @(\"16+19=35\", \"50-6=44\"),\n    @(\"85-84=1\", \"64+22=86\"),\n    @(\"84-23=61\", \"25+60=85\"),\n    @(\"18+60=78\", \"29+56=85\"),\n    @(\"4+62=66\", \"64-26=38\"),\n    @(\"99-56=43\", \"40+18=58\"),\n    @(\"8+56=64\", \"66-2=64\"),\n    @(\"47-42=5\", \"90-24=66\"),\n    @(\"89-17=72\", \"85-19=66\"),\n    @(\"2+16=18\", \"29+8=37\"),\n    @(\"34+26=60\", \"19-11=8\"),\n    @(\"20+22=42\", \"48+0=48\"),\n    @(\"26+57=83\", \"40+9=49\"),\n    @(\"39-14=25\", \"98-11=87\"),\n    @(\"98-59=39\", \"19+36=55\"),\n    @(\"78-12=66\", \"80-25=55\"),\n    @(\"6+66=72\", \"61+11=72\"),\n    @(\"58+21=79\", \"98-81=17\"),\n    @(\"83-23=60\", \"13+71=84\"),\n    @(\"78-56=22\", \"77-21=56\"),\n    @(\"62-20=42\", \"24+33=57\"),\n    @(\"80+4=84\", \"24-0=24\"),\n    @(\"35-0=35\", \"11-3=8\"),\n    @(\"59-51=8\", \"0+52=52\"),\n    @(\"92-41=51\", \"66-65=1\"),\n    @(\"56-4=52\", \"94-48=46\"),\n    @(\"22+72=94\", \"33+52=85\"),\n    @(\"59-4=55\", \"73+2=75\"),\n    @(\"98-86=12\", \"95+3=98\"),\n    @(\"76-22=54\", \"10-7=3\"),\n    @(\"94+3=97\", \"43+50=93\"),\n    @(\"32+38=70\", \"59-15=44\"),\n    @(\"94-23=71\", \"21+29=50\"),\n    @(\"74-23=51\", \"42+30=72\"),\n    @(\"81-20=61\", \"21+39=60\"),\n    @(\"81-4=77\", \"31+23=54\"),\n    @(\"57-17=40\", \"89-85=4\"),\n    @(\"55-18=37\", \"44+54=98\"),\n    @(\"15+80=95\", \"39-2=37\"),\n    @(\"28-6=22\", \"86-22=64\"),\n    @(\"38+47=85\", \"37+29=66\"),\n    @(\"32+22=54\", \"44+46=90\"),\n    @(\"23+75=98\", \"25-19=6\"),\n    @(\"12+16=28\", \"82-63=19\"),\n    @(\"28+29=57\", \"84-74=10\"),\n    @(\"65-57=8\", \"28+45=73\"),\n    @(\"54-10=44\", \"6+14=20\"),\n    @(\"2+42=44\", \"50+21=71\"),\n    @(\"37+30=67\", \"26+1=27\"),\n    @(\"37+48=85\", \"24+48=72\"),\n    @(\"16+43=59\", \"92-88=4\"),\n    @(\"95-71=24\", \"84-81=3\"),\n    @(\"44+14=58\", \"56-16=40\"),\n    @(\"18-5=13\", \"61-40=21\"),\n    @(\"27-2=25\", \"9+1=10\"),\n    @(\"60-38=22\", \"12+58=70\"),\n    @(\"89-47=42\", \"90-46=44\"),\n    @(\"22-15=7\", \"52-46=6\"),\n    @(\"91-0=91\", \"67-2=65\"),\n    @(\"78-36=42\", \"78-64=14\"),\n    @(\"61+28=89\", \"52+15=67\"),\n    @(\"63-29=34\", \"40+44=84\"),\n    @(\"70-43=27\", \"37-22=15\"),\n    @(\"79-4=75\", \"89-80=9\"),\n    @(\"9+57=66\", \"7+78=85\"),\n    @(\"29+13=42\", \"0+96=96\"),\n    @(\"80-20=60\", \"28+61=89\"),\n    @(\"38-17=21\", \"26+25=51\"),\n    @(\"81-43=38\", \"89-33=56\"),\n    @(\"8+37=45\", \"16+12=28\"),\n    @(\"19+30=49\", \"77+3=80\"),\n    @(\"81+14=95\", \"7+11=18\"),\n    @(\"79-67=12\", \"5+40=45\"),\n    @(\"61-0=61\", \"16+53=69\"),\n    @(\"3+16=19\", \"71-69=2\"),\n    @(\"35+48=83\", \"34+49=83\"),\n    @(\"44+53=97\", \"19+51=70\"),\n    @(\"8+11=19\", \"41-15=26\"),\n    @(\"33-2=31\", \"20+19=39\"),\n    @(\"82+6=88\", \"21+31=52\")\n)\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif (($rowCount * $colCount) -ne $pairs.Length) {\n    throw \"Cell count ($($rowCount * $colCount)) does not match expected pair count ($($pairs.Length)).\"\n}\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cellRange = $cell.Range\n        # Cell.Range.Text includes the trailing end-of-cell marker (CR + BEL,\n        # chars 13/7); strip it so we can compare/replace the visible text only.\n        $rawText = $cellRange.Text\n        $currentText = $rawText.TrimEnd([char]13, [char]7)\n\n        $pair = $pairs[$i]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n\n        if ($currentText -ne $oldText) {\n            throw \"Cell $i (row $r, col $c) text mismatch: expected '$oldText' but found '$currentText'.\"\n        }\n\n        # Assigning to Range.Text replaces the run's text in place, preserving\n        # the existing run/paragraph formatting (font, size, alignment, ...)\n        # instead of wiping it out, matching how the diff only touches <w:t>.\n        $cellRange.Text = $newText\n\n        $i++\n    }\n}\n"}
